# Update "想去人数" (F column) figures across the four worksheets to match
# the latest scrape output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 784
$ws.Range("F3").Value = 2835
$ws.Range("F4").Value = 1342
$ws.Range("F11").Value = 11802
$ws.Range("F12").Value = 6688
$ws.Range("F20").Value = 93
$ws.Range("F21").Value = 284
$ws.Range("F22").Value = 933
$ws.Range("F23").Value = 3659
$ws.Range("F25").Value = 990
$ws.Range("F27").Value = 177
$ws.Range("F31").Value = 271
$ws.Range("F32").Value = 313
$ws.Range("F33").Value = 5042
$ws.Range("F35").Value = 1252
$ws.Range("F36").Value = 243
$ws.Range("F37").Value = 571
$ws.Range("F38").Value = 208
$ws.Range("F39").Value = 548

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 3694

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9078
$ws.Range("F4").Value = 1842

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9078
$ws.Range("F4").Value = 1842
$ws.Range("F5").Value = 784
$ws.Range("F6").Value = 2835
$ws.Range("F10").Value = 1342
$ws.Range("F17").Value = 11803
$ws.Range("F18").Value = 3694
$ws.Range("F19").Value = 6688
$ws.Range("F27").Value = 93
$ws.Range("F28").Value = 933
$ws.Range("F29").Value = 3659
$ws.Range("F31").Value = 990
$ws.Range("F32").Value = 177
$ws.Range("F38").Value = 313
$ws.Range("F39").Value = 1252
$ws.Range("F40").Value = 243
$ws.Range("F41").Value = 208
$ws.Range("F42").Value = 548

$wb.Save()
